# Updates cryptos list data (price and volume-change columns) per commit
# "Updated cryptos list on Sun Sep 10 19:41:33 UTC 2023 with GitHub Actions".
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cells in column D whose new values look like plain numbers (e.g. "213.46").
# Excel auto-converts such text to a numeric value on assignment, but the
# source file stores them as plain text (inline strings). Temporarily mark
# these cells as Text ("@") before assigning, then restore the original
# "Normal" style so no visible formatting change is left behind.
$textForceCells = @("D5", "D6", "D10", "D17", "D20", "D22", "D23", "D25", "D28", "D29", "D31", "D32", "D37", "D39", "D41", "D46", "D51")
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).NumberFormat = "@"
}

$values = @{
    "D2" = "25.898.17"
    "E2" = "  -0.37%  "
    "D3" = "1.621.42"
    "E3" = "  -1.22%  "
    "E4" = "  +0.30%  "
    "D5" = "213.46"
    "E5" = "  -1.14%  "
    "D6" = "0.501"
    "E6" = "  -1.12%  "
    "E7" = "  +0.30%  "
    "E8" = "  -2.54%  "
    "E9" = "  -3.80%  "
    "D10" = "18.17"
    "E10" = "  -6.87%  "
    "E11" = "  -1.24%  "
    "D12" = "1.848.42"
    "E12" = "  -1.12%  "
    "D13" = "1.623.07"
    "E13" = "  -1.28%  "
    "E14" = "  -2.42%  "
    "E15" = "  -3.82%  "
    "D16" = "25.892.67"
    "E16" = "  -0.50%  "
    "D17" = "60.99"
    "E17" = "  -3.75%  "
    "E18" = "  -4.16%  "
    "E19" = "  +0.25%  "
    "D20" = "191.90"
    "E20" = "  -1.13%  "
    "E21" = "  -3.30%  "
    "D22" = "9.53"
    "E22" = "  -3.91%  "
    "D23" = "6.05"
    "E23" = "  -2.45%  "
    "E24" = "  +0.31%  "
    "D25" = "143.64"
    "E25" = "  +0.48%  "
    "E26" = "  +0.35%  "
    "E27" = "  -3.33%  "
    "D28" = "6.69"
    "E28" = "  -2.61%  "
    "D29" = "15.10"
    "E29" = "  -2.69%  "
    "E30" = "  -1.47%  "
    "D31" = "0.0481"
    "E31" = "  -2.89%  "
    "D32" = "3.10"
    "E32" = "  -4.80%  "
    "E33" = "  -5.95%  "
    "E34" = "  -3.20%  "
    "E35" = "  -2.38%  "
    "D36" = "1.115.11"
    "E36" = "  -1.26%  "
    "D37" = "0.842"
    "E37" = "  -6.87%  "
    "E38" = "  -1.31%  "
    "D39" = "0.516"
    "E39" = "  -4.37%  "
    "E40" = "  -2.68%  "
    "D41" = "97.87"
    "E41" = "  -1.11%  "
    "E42" = "  -4.16%  "
    "D43" = "1.758.82"
    "E43" = "  -1.08%  "
    "E44" = "  -5.80%  "
    "E45" = "  -1.67%  "
    "D46" = "0.0530"
    "E46" = "  +1.55%  "
    "E47" = "  -4.17%  "
    "E48" = "  -2.16%  "
    "E49" = "  -0.32%  "
    "E50" = "  +0.43%  "
    "D51" = "7.45"
    "E51" = "  -3.66%  "
}

foreach ($cellRef in $values.Keys) {
    $ws.Range($cellRef).Value = $values[$cellRef]
}

# Restore default ("Normal") style on the cells we temporarily reformatted,
# so only the cell *values* changed -- not their formatting.
foreach ($cellRef in $textForceCells) {
    $ws.Range($cellRef).Style = "Normal"
}
